# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
#
# The underlying match rows got re-keyed (ids changed / rows re-ordered),
# which manifests as several small groups of rows exchanging their data
# (everything except the row-index column A, which stays put).
#
# For every group below, row[i] ends up holding the data that row[i+1]
# held before the edit (and the last row in the group wraps around to
# the data that the first row held before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows {
    # Positional parameter: ordered list of row numbers; row[i] receives
    # the pre-edit contents of row[i+1] (wrapping around).
    param($Rows)

    $firstCol = "B"
    $lastCol  = "AD"

    # Snapshot every row's current (pre-edit) B:AD contents up front so
    # that later writes don't clobber data we still need to read.
    $snapshots = @{}
    foreach ($r in $Rows) {
        $rng = $ws.Range("$firstCol$r`:$lastCol$r")
        $snapshots[$r] = $rng.Value2
    }

    $count = $Rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $Rows[$i]
        $sourceRow = $Rows[($i + 1) % $count]
        $destRng = $ws.Range("$firstCol$targetRow`:$lastCol$targetRow")
        $destRng.Value2 = $snapshots[$sourceRow]
    }
}

# Rows 93-96 rotate as a 4-cycle: 93<-96, 96<-95, 95<-94, 94<-93
Rotate-Rows @(93, 96, 95, 94)

# Remaining groups are simple two-row swaps
Rotate-Rows @(100, 101)
Rotate-Rows @(102, 103)
Rotate-Rows @(114, 115)
Rotate-Rows @(162, 163)
Rotate-Rows @(173, 174)
Rotate-Rows @(205, 206)
